# "divide variable beta 1"
# Binary Instructions workbook: change row 8's instruction, push its old
# content down into a newly-populated row 9, and populate row 10 with a
# third variant. Columns R:S and X:AN recompute automatically via the
# sheet's existing shared formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Row 10: brand-new instruction (cells were blank, formatting already
#     matches what the final file needs, so plain value writes are enough).
$ws.Range("B10").Value = "00110"
$ws.Range("C10").Value = "000"
$ws.Range("D10").Value = "111"
$ws.Range("E10").Value = "0"
$ws.Range("F10").Value = "0"
$ws.Range("G10").Value = "1"
$ws.Range("H10").Value = "000"
$ws.Range("I10").Value = "000"
$ws.Range("J10").Value = "0"
$ws.Range("K10").Value = "0"
$ws.Range("L10").Value = "0"
$ws.Range("M10").Value = "000"
$ws.Range("N10").Value = "000"
$ws.Range("O10").Value = "0"
$ws.Range("P10").Value = "0"
$ws.Range("Q10").Value = "0"

# --- Row 8: only Instr/B and the B-register Alt/Write bit change.
$ws.Range("B8").Value = "01001"
$ws.Range("I8").Value = "001"
$ws.Range("J8").Value = "1"

# Writing text through COM always lands on the "quote-prefixed" twin of the
# cell's style; restore the original (non-quote-prefixed) look by pulling
# formatting back from the untouched row above, which still carries it.
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial($fmt)
$ws.Range("I7").Copy()
$ws.Range("I8").PasteSpecial($fmt)
$ws.Range("J7").Copy()
$ws.Range("J8").PasteSpecial($fmt)
$excel.CutCopyMode = $false

# --- Row 9: becomes what row 8 used to hold (same values). Write the
#     values first, then re-paste row 8's formats over them last, since
#     every text write re-quote-prefixes whichever cell it touches.
$ws.Range("B9").Value = "11000"
$ws.Range("C9").Value = "000"
$ws.Range("D9").Value = "000"
$ws.Range("E9").Value = "0"
$ws.Range("F9").Value = "0"
$ws.Range("G9").Value = "0"
$ws.Range("H9").Value = "000"
$ws.Range("I9").Value = "000"
$ws.Range("J9").Value = "0"
$ws.Range("K9").Value = "0"
$ws.Range("L9").Value = "0"
$ws.Range("M9").Value = "000"
$ws.Range("N9").Value = "001"
$ws.Range("O9").Value = "0"
$ws.Range("P9").Value = "1"
$ws.Range("Q9").Value = "0"

$ws.Range("B8:Q8").Copy()
$ws.Range("B9:Q9").PasteSpecial($fmt)
$excel.CutCopyMode = $false

# --- Selection moves to the block of recomputed hex results.
$ws.Range("S2:S10").Select()
